$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fifa_world_cup_2018_matches")

# Update status (column D) from "open" to "completed" for rows 3-9
# and update goals (columns G/H) for the same rows.
$ws.Range("D3").Value = "completed"
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = "completed"
$ws.Range("H4").Value = 1

$ws.Range("D5").Value = "completed"
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 3

$ws.Range("D6").Value = "completed"
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 1

$ws.Range("D7").Value = "completed"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1

$ws.Range("D8").Value = "completed"
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = "completed"
$ws.Range("G9").Value = 2

# Update the active selection on the active sheet's bottom-right pane
$ws.Range("G5").Select()
